$d = $word.ActiveDocument
$sec = $d.Sections(1)

# Rename the inline picture identified by its docPr/Shape ID, found by
# scanning the given header/footer range's InlineShapes collection.
# Word's InlineShape object has no direct .Name setter, so we briefly
# convert the picture to a floating Shape (which does expose .Name --
# this maps to the <wp:docPr name="..."> attribute), rename it, then
# convert it straight back to an inline picture so the layout/anchoring
# is left exactly as it was.
function Rename-InlineLogoById($range, $targetId, $newName) {
    $count = $range.InlineShapes.Count
    for ($i = 1; $i -le $count; $i++) {
        $ishp = $range.InlineShapes($i)
        $shp = $ishp.ConvertToShape()
        if ($shp.ID -eq $targetId) {
            $shp.Name = $newName
            $shp.ConvertToInlineShape() | Out-Null
            return
        }
        $shp.ConvertToInlineShape() | Out-Null
    }
}

# Pearson logo pictures (footers): image1.png -> image2.png
Rename-InlineLogoById $sec.Footers(1).Range 4 "image2.png"
Rename-InlineLogoById $sec.Footers(2).Range 2 "image2.png"

# BTec logo pictures (headers): image2.jpg -> image1.jpg
Rename-InlineLogoById $sec.Headers(1).Range 3 "image1.jpg"
Rename-InlineLogoById $sec.Headers(2).Range 1 "image1.jpg"

Write-Host "Renamed inline picture labels in headers/footers."
